$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('outros', -0.08671325966487606, 3.987593189047669, -6.650134709794301, -6.627930525235978),
    @('vida afetiva', 1.265342662830762, 1.957417104250754, -2.643620611784112, -6.33376009333355),
    @('recreação, lazer e entretenimento', 0.8663547429634153, 2.362193171803985, -4.557644810472275, -5.030313185036145),
    @('cultural', 0.2147805530085035, 1.471414819249782, -1.623195302797334, -3.742836356467645),
    @('linguagem e representação', 1.44543401503942, 1.153204681192293, -3.308428478809178, -3.289841908814284),
    @('sociedade e consumo', 0.3060084338911121, 1.398620750262873, -2.528762197866694, -2.759868412839082),
    @('dinâmica urbana', 2.476480745930567, -1.004457874634521, 1.336979290858656, -2.603135590561373),
    @('estruturas sociais e econômicas', 0.5889691512491492, 1.355391004008771, -3.08186292015781, -2.590608708207778),
    @('saúde', 1.26506595475886, 0.1416705975399979, -0.5609731936916506, -2.292379325518055),
    @('família', -3.232767059716722, 0.08952967999095393, 7.904587170393286, -2.159355282519024),
    @('gênero e sexualidade', -1.794620843912653, 2.510081912336033, -2.961694015273024, -2.131263919694551),
    @('vida rural, vida no campo', -0.2328921378020009, 0.906188594987037, -0.5182084830224819, -2.102033170279678),
    @('nacionalidade e nacionalismo', 1.673256289713538, 0.3433577954168633, -2.363784615723356, -1.931424510879601),
    @('trabalho e ofício', -0.09865352015420724, 0.3180583477723646, 0.8462121018510294, -1.737545538037684),
    @('distopia, ficção científica e futuros imaginados', 4.032814949354573, -1.228380860505598, -2.010011465820939, -1.413073224251723),
    @('política', -2.114553349931825, 1.678864341440947, -0.3226874853552322, -1.371564699918195),
    @('poesia e ensaio', 0.6395591056517784, 0.2322421338701065, -0.6312503708577855, -1.35240822690325),
    @('desinformação, populismo e polarização', 0.003459243809323817, 0.7251421314193658, -1.234409691415849, -1.216243992744902),
    @('violência', 0.336598511868042, 0.1823141054070368, -0.1115436180749076, -1.136057071210128),
    @('memória e patrimônio', 0.03751443314315746, 0.635689554350033, -1.299344419500882, -0.9130596725497483),
    @('tecnologia, inovação e sociedade', -0.30671195187988, 0.1506116250340506, 0.9754483592135256, -0.9109487337521832),
    @('questão agrária e territorial', -0.8643779686337749, 0.283463408537552, 1.530871392945483, -0.8689893294072423),
    @('sistema penitenciário', -1.174053302241971, 0.7265157921018225, 0.2601699362852061, -0.5080601675090328),
    @('educação e socialização', 0.6530256318040129, -0.1689308550567603, -0.3447845900911859, -0.3106596862147576),
    @('pandemia', 0.9884617807957636, -1.562865413229943, 3.439107497960149, -0.03114723011075469),
    @('mobilidade', -0.54462359778988, -0.7979905197469427, 3.603709064852528, 0.1017489805363181),
    @('sonoridade e paisagem sonora', 0.657858089643153, -0.3364265371722031, -0.4823759363756056, 0.3838268459731174),
    @('amizade', -1.097550062859834, 0.1673253246593815, 0.7414466546828354, 0.7532989000489488),
    @('territorialidade e colonialismo', 0.6068606657190214, -1.30107740935588, 2.456754333556304, 0.7876087443358792),
    @('alimentação e tratamentos tradicionais', -0.6789962442763704, 0.274947376382655, -0.6285055794253361, 0.9780637733247727),
    @('raça e racismo', -1.644155044596602, -0.2500033643228159, 2.638791089415573, 1.285627346229134),
    @('pessoas com deficiência', 1.464165188471019, -1.647121534755042, 1.396492849795139, 1.408992354145778),
    @('povos originários e comunidades tradicionais', 0.306908989172916, -0.579624721024082, -0.4204713586112711, 1.802214587302508),
    @('maternidade', -0.1801812774922617, -0.5975923142142956, 0.3656195921450998, 1.989339573138265),
    @('saúde mental', 0.8996540907221953, -1.674291009308457, 1.767110042445613, 2.189515870340706),
    @('religião, espiritualidade e cosmologias', -0.4150204948944516, -0.6327688170384438, 0.7075669737000131, 2.206248097314158),
    @('ambiental', -0.2141817334953964, -1.526066821794413, 3.269868792753759, 2.274042211884813),
    @('sonho e fantasia', 0.8394204789531036, -1.065257184346216, -0.2799859896477954, 2.300367439809822),
    @('violências e preconceitos de gênero', -1.140965811655277, 0.162513281035833, -0.8378317108637048, 2.435687958705589),
    @('saudade, luto ou perda', 1.887103886674322, -1.4701655795871, -1.248903672666485, 2.671144386748706),
    @('moradia e habitação', -1.518569550488702, -2.958953561051473, 10.05650120851277, 2.753306088439786),
    @('crises e desastres ambientais e sociais', -0.2596433922561583, -1.237627908048992, 1.806425111619174, 2.853952510141601),
    @('direitos humanos', -0.5889832084098315, -0.6418142331413879, 0.2351070225803011, 3.037686359834469),
    @('corpo, performance e expressão', -1.951319644673254, 0.5381099314946762, -1.980925567242882, 3.837048410058943),
    @('biografia', -0.1307284400838844, -1.491706855660136, -1.723967114376707, 7.011799776735748),
    @('solidão', -1.463650198192873, -2.930504827727594, 5.313745162880253, 7.312712619449605),
    @('vida cotidiana', -2.372092927874691, -1.625860500957157, 2.08835681934992, 7.850907350069113),
    @('mulher, feminino e feminismo', -1.345888291968431, -3.133016099018707, 4.933491104990471, 8.156716340786859),
    @('resistência e luta', -2.142224056928094, -3.001834732383529, 5.454175556631077, 8.687603048609532),
    @('reflexão', -3.086926411141425, -4.196028148474826, 4.609832840981095, 15.34170193867449),
    @('arte', -1.891097508517258, -5.777095969039737, 5.200954688781274, 17.83991940669465)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $r = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}

Write-Output "Done updating rows 2-52"
